$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Header: rename "Estimativa" -> "Estimativa em horas"
$ws.Range("D1").Value = "Estimativa em horas"

# Update use-case list (column B), re-ordering UC009 to the bottom of the
# list and renaming UC006 to "Consultar informações", then fill in the new
# "hours" estimates (column D) for every use case.
$ws.Range("B2").Value = "UC004 - Manter Cotação"
$ws.Range("D2").Value = 79

$ws.Range("B3").Value = "UC006 - Consultar informações"
$ws.Range("D3").Value = 39

$ws.Range("B4").Value = "UC009 - Manter Apólice"
$ws.Range("D4").Value = 39

$ws.Range("B5").Value = "UC001 - Manter Segurado"
$ws.Range("D5").Value = 79

$ws.Range("B6").Value = "UC002 - Manter Objeto de Seguro"
$ws.Range("D6").Value = 79

$ws.Range("B7").Value = "UC003 - Manter aditamento/clausula"
$ws.Range("D7").Value = 79

$ws.Range("B8").Value = "UC005 - Manter ocorrência"
$ws.Range("D8").Value = 79

$ws.Range("B9").Value = "UC007 - Manter regra"
$ws.Range("D9").Value = 79

$ws.Range("B10").Value = "UC008 - Manter sinistro"
$ws.Range("D10").Value = 39

$ws.Range("B11").Value = "UC010 - Gerar renovação apolice"
$ws.Range("D11").Value = 39

# Update the explanatory note text in the merged A13:D21 block.
$ws.Range("A13").Value = "A tecnica de estimativa utilizada foi UCP, devido a facilidade de estimativa do software com casos de uso já desenvolvidos;"

# New underlined, empty cell next to the note block (style/font addition).
$ws.Range("E13").Font.Underline = $true

# Widen column D to fit the new, longer header text.
$ws.Columns.Item(4).ColumnWidth = 19.109375

# Matches the cursor position left behind in the saved file.
$null = $ws.Range("H27").Select()

$wb.Save()
